# ValueSet-differential-unspecified-blood-vs.xlsx update
# - bump Version, Status, Date
# - replace Contact row text, add a second Contact row (Bob Milius)
# - insert a new "Jurisdiction" row after Contact rows
# (sheet2 "Include from LOINC" content is untouched - its shared-string
#  indices merely shift because of the sharedStrings.xml edits above)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Version 0.1.6 -> 0.1.7
$ws.Range("B3").Value = "0.1.7"

# Status active -> draft
$ws.Range("B6").Value = "draft"

# Date
$ws.Range("B8").Value = "2024-08-27T12:23:18-05:00"

# Contact (row 10) gets the publisher-with-URL text
$ws.Range("B10").Value = "The Medical College of Wisconsin, Inc. and the National Marrow Donor Program (http://www.cibmtr.org)"

# Contact (row 11) gets the named contact
$ws.Range("B11").Value = "Bob Milius (bmilius@nmdp.org)"

# Insert a new row 12 for "Jurisdiction" (empty value), pushing
# Description/Purpose/Copyright/Immutable down by one row each.
$ws.Rows.Item(12).Insert()

# Restore the bordered/wrapped row style on the newly inserted row
$ws.Range("A11:B11").Copy()
$ws.Range("A12:B12").PasteSpecial(-4122)

$ws.Range("A12").Value = "Jurisdiction"
$ws.Range("B12").Value = ""

Write-Host "edit complete"
